$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns at D:E, shifting existing D:K data to F:M
$ws.Range("D:E").Insert()

# Copy formatting (styles/number formats) from column F into the new D:E columns
$ws.Range("F7:F102").Copy() | Out-Null
$ws.Range("D7:E102").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Populate the new D:E columns with the new quarter data (Dec-2018 / Sep-2018 columns)
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("D8").Value = 45200
$ws.Range("E8").Value = 37100
$ws.Range("D9").Value = 15000
$ws.Range("E9").Value = 12100
$ws.Range("D10").Value = 30200
$ws.Range("E10").Value = 25000
$ws.Range("D12").Value = 5700
$ws.Range("E12").Value = 5400
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = 10000
$ws.Range("E14").Value = 2500
$ws.Range("D15").Value = 4700
$ws.Range("E15").Value = 3600
$ws.Range("D17").Value = 50000
$ws.Range("E17").Value = 36800
$ws.Range("D18").Value = -4800
$ws.Range("E18").Value = 300
$ws.Range("D20").Value = -5400
$ws.Range("E20").Value = -3800
$ws.Range("D21").Value = -3400
$ws.Range("E21").Value = 1900
$ws.Range("D22").Value = 0
$ws.Range("E22").Value = 0
$ws.Range("D23").Value = -10100
$ws.Range("E23").Value = -3500
$ws.Range("D24").Value = -11900
$ws.Range("E24").Value = 700
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = 1800
$ws.Range("E26").Value = -4300
$ws.Range("D27").Value = 1800
$ws.Range("E27").Value = -4300
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = "NA"
$ws.Range("E29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = 5400
$ws.Range("E32").Value = 3800
$ws.Range("D33").Value = 1800
$ws.Range("E33").Value = -4300
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = 1800
$ws.Range("E35").Value = -4300
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("D41").Value = 16700
$ws.Range("E41").Value = 16100
$ws.Range("D42").Value = 0
$ws.Range("E42").Value = 0
$ws.Range("D43").Value = 44500
$ws.Range("E43").Value = 26400
$ws.Range("D44").Value = 0
$ws.Range("E44").Value = 0
$ws.Range("D45").Value = 6000
$ws.Range("E45").Value = 6300
$ws.Range("D46").Value = 67300
$ws.Range("E46").Value = 48800
$ws.Range("D47").Value = 1600
$ws.Range("E47").Value = 1600
$ws.Range("D48").Value = 2800
$ws.Range("E48").Value = 2200
$ws.Range("D49").Value = 404900
$ws.Range("E49").Value = 269200
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 6600
$ws.Range("E52").Value = 5600
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 483200
$ws.Range("E54").Value = 327500
$ws.Range("D57").Value = 3500
$ws.Range("E57").Value = 4100
$ws.Range("D58").Value = 6000
$ws.Range("E58").Value = 4300
$ws.Range("D59").Value = 98100
$ws.Range("E59").Value = 67400
$ws.Range("D60").Value = 107600
$ws.Range("E60").Value = 75900
$ws.Range("D61").Value = 273700
$ws.Range("E61").Value = 153900
$ws.Range("D62").Value = 14500
$ws.Range("E62").Value = 8400
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 395900
$ws.Range("E66").Value = 238300
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = -85700
$ws.Range("E72").Value = -87500
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 87300
$ws.Range("E76").Value = 89200
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("D81").Value = 1800
$ws.Range("E81").Value = -4300
$ws.Range("D83").Value = 6700
$ws.Range("E83").Value = 5400
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = 2600
$ws.Range("E89").Value = 1100
$ws.Range("D91").Value = -300
$ws.Range("E91").Value = -100
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = -113200
$ws.Range("E94").Value = -2600
$ws.Range("D96").Value = 0
$ws.Range("E96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = 112400
$ws.Range("E100").Value = -1700
$ws.Range("D101").Value = -1100
$ws.Range("E101").Value = 300
$ws.Range("D102").Value = 600
$ws.Range("E102").Value = -2900
